$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11; this shifts the existing rows 11-66
# (and all their data/formatting) down to rows 12-67, matching the
# target diff where every row from the old 11 downward is pushed down
# by one and a fresh weekly observation is inserted at row 11.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new weekly price record.
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C11").Value = "Arica y Parinacota"
$ws.Range("D11").Value = 44592
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 100112040
$ws.Range("G11").Value = "Cilantro"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 1500
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 1750
$ws.Range("N11").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 875
$ws.Range("Q11").Value = 2
$ws.Range("R11").Value = "Hortaliza"
